# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.702.47"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "2.047.28"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  +0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "227.25"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.608"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.71%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "58.90"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.57%  "
$ws.Range("E8").Value = "  +0.05%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.375"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.08%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0829"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "2.352.44"
$ws.Range("E12").Value = "  +0.95%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "14.32"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.78%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "21.16"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.27%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.47"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +5.64%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.757"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").Value = "2.048.03"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").Value = "37.646.28"
$ws.Range("E18").Value = "  -0.41%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "5.98"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.50%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "69.28"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").Value = "0.0₃0827"
$ws.Range("E21").Value = "  +0.35%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "222.24"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("E24").Value = "  +0.05%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.27"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.58%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "169.08"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.67%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.25"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").Value = "  +0.18%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "18.70"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.01%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.28"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  -1.19%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.37"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +12.80%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.35"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.82%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.54"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.32%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0599"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.87%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "6.39"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.30%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.35"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +3.64%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.42"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +4.03%  "
$ws.Range("E39").Value = "  +0.10%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "18.19"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +7.85%  "
$ws.Range("D41").Value = "1.535.81"
$ws.Range("E41").Value = "  +0.73%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "99.27"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.88%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.0214"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.96%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.83"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0892"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "4.09"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.95"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.89%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "7.08"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").Value = "2.240.91"
$ws.Range("E51").Value = "  +0.97%  "
